$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was "M") -> becomes "B"
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9448275862068966
$ws.Range("C2").Value = 0.958041958041958
$ws.Range("D2").Value = 0.9513888888888888
$ws.Range("E2").Value = 143

# Row 3 (was "B") -> becomes "M"
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.927710843373494
$ws.Range("C3").Value = 0.9058823529411765
$ws.Range("D3").Value = 0.9166666666666666
$ws.Range("E3").Value = 85

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9385964912280702
$ws.Range("C4").Value = 0.9385964912280702
$ws.Range("D4").Value = 0.9385964912280702
$ws.Range("E4").Value = 0.9385964912280702

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9362692147901952
$ws.Range("C5").Value = 0.9319621554915672
$ws.Range("D5").Value = 0.9340277777777777

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9384463443611106
$ws.Range("C6").Value = 0.9385964912280702
$ws.Range("D6").Value = 0.9384442007797271
